$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 61
$ws.Range("F2").Value = 43
$ws.Range("H2").Value = 55

# Row 9
$ws.Range("E9").Value = 31

# Row 15
$ws.Range("E15").Value = 178

# Row 18
$ws.Range("E18").Value = 134

# Row 19
$ws.Range("E19").Value = 70

# Row 38
$ws.Range("E38").Value = 86

# Row 40
$ws.Range("E40").Value = 28

# Row 46
$ws.Range("E46").Value = 31

# Row 47
$ws.Range("F47").Value = 41
$ws.Range("H47").Value = 51

# Row 49
$ws.Range("E49").Value = 78

# Row 63
$ws.Range("E63").Value = 46

# Row 70
$ws.Range("E70").Value = 50

# Row 77
$ws.Range("E77").Value = 65

# Row 79
$ws.Range("E79").Value = 45
$ws.Range("F79").Value = 22
$ws.Range("H79").Value = 33
